$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Table layout (Nationality x Period):
#   Row 1  : header (periods)
#   Row 2  : "Great Britain" section header
#   Row 3  : "  STDT"  data row for Great Britain
#   Row 4  : "  Sample" data row for Great Britain
#   Row 5  : "Netherlands" section header
#   Row 6  : "  STDT"  data row for Netherlands
#   Row 7  : "  Sample" data row for Netherlands
#   Row 8  : "France" section header
#   Row 9  : "  STDT"  data row for France
#   Row 10 : "  Sample" data row for France
#   Row 11 : "Total" section header
#   Row 12 : "  STDT"  data row for Total
#   Row 13 : "  Sample" data row for Total
#
# Column 1 holds the row label ("  STDT" / "  Sample"); columns 2-11 hold the
# count/fraction pairs for the five period columns (1750-1762, 1763-1778,
# 1778-1783, 1784-1795, Total).

# 1) "STDT" -> "TSTD" label, once per nationality block (Great Britain,
#    Netherlands, France, Total)
foreach ($r in 3, 6, 9, 12) {
    $cell = $tbl.Cell($r, 1)
    $cell.Range.Text = "  TSTD"
}

# 2) Great Britain / Sample row (row 4): 123 -> 122 (1763-1778 count),
#    208 -> 207 (Total count), 0.56 -> 0.55 (Total fraction)
$tbl.Cell(4, 4).Range.Text = "122"
$tbl.Cell(4, 10).Range.Text = "207"
$tbl.Cell(4, 11).Range.Text = "0.55"

# 3) France / Sample row (row 10): 0.11 -> 0.12 (1763-1778 fraction)
$tbl.Cell(10, 5).Range.Text = "0.12"

# 4) Total / Sample row (row 13): 216 -> 215 (1763-1778 count),
#    374 -> 373 (Total count)
$tbl.Cell(13, 4).Range.Text = "215"
$tbl.Cell(13, 10).Range.Text = "373"
